$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H column (Ligand total expression value) - fix floating point representation
$ws.Range("H2").Value = 0.09138
$ws.Range("H3").Value = 0.09138
$ws.Range("H4").Value = 0.09138

# M2 / N2 (Receptor average / total expression value) - the root data change
$ws.Range("M2").Value = 3.636278
$ws.Range("N2").Value = 10.908834

# Recompute derived specificity / edge weight columns from M/N (and G/H, I/J)
$G2 = $ws.Range("G2").Value()
$G3 = $ws.Range("G3").Value()
$G4 = $ws.Range("G4").Value()

$H2 = $ws.Range("H2").Value()
$H3 = $ws.Range("H3").Value()
$H4 = $ws.Range("H4").Value()

$I2 = $ws.Range("I2").Value()
$I3 = $ws.Range("I3").Value()
$I4 = $ws.Range("I4").Value()

$J2 = $ws.Range("J2").Value()
$J3 = $ws.Range("J3").Value()
$J4 = $ws.Range("J4").Value()

$M2 = $ws.Range("M2").Value()
$M3 = $ws.Range("M3").Value()
$M4 = $ws.Range("M4").Value()

$N2 = $ws.Range("N2").Value()
$N3 = $ws.Range("N3").Value()
$N4 = $ws.Range("N4").Value()

$sumM = $M2 + $M3 + $M4
$sumN = $N2 + $N3 + $N4

$O2 = $M2 / $sumM
$O3 = $M3 / $sumM
$O4 = $M4 / $sumM

$P2 = $N2 / $sumN
$P3 = $N3 / $sumN
$P4 = $N4 / $sumN

$ws.Range("O2").Value = $O2
$ws.Range("O3").Value = $O3
$ws.Range("O4").Value = $O4

$ws.Range("P2").Value = $P2
$ws.Range("P3").Value = $P3
$ws.Range("P4").Value = $P4

$ws.Range("Q2").Value = $G2 * $M2
$ws.Range("Q3").Value = $G3 * $M3
$ws.Range("Q4").Value = $G4 * $M4

$ws.Range("R2").Value = $H2 * $N2
$ws.Range("R3").Value = $H3 * $N3
$ws.Range("R4").Value = $H4 * $N4

$ws.Range("S2").Value = $I2 * $O2
$ws.Range("S3").Value = $I3 * $O3
$ws.Range("S4").Value = $I4 * $O4

$ws.Range("T2").Value = $J2 * $P2
$ws.Range("T3").Value = $J3 * $P3
$ws.Range("T4").Value = $J4 * $P4
